# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-19 09:29:12
#
# The "Recorded By" column (G) lists the users who touched a session
# record. Two specific combinations had their listed order corrected:
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#   "system, backup@backdoor.com, System"  -> "backup@backdoor.com, System, system"
# All other values in column G are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    }
    elseif ($text -eq "system, backup@backdoor.com, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    }
}
